# Rename header labels on the existing sheets
$wb = $excel.ActiveWorkbook

$weekly = $wb.Worksheets.Item(1)
$weekly.Range("B1").Value = "Weekly_PO_Qty"

$monthly = $wb.Worksheets.Item(2)
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet right after "Monthly Trend"
$poForecast = $wb.Worksheets.Add($null, $monthly)
$poForecast.Name = "PO Forecast"

# Copy header-row formatting (bold + border + centered) from an existing
# header cell instead of re-building the style by hand, so the new sheet
# reuses the workbook's existing header style.
$weekly.Range("B1").Copy()
$poForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-column formatting from the existing date column too.
$weekly.Range("A2").Copy()
$poForecast.Range("A2:A11").PasteSpecial(-4122)

# Header labels
$poForecast.Range("A1").Value = "ds"
$poForecast.Range("B1").Value = "PO_Forecast"
$poForecast.Range("C1").Value = "yhat_lower"
$poForecast.Range("D1").Value = "yhat_upper"

# Forecast rows
$poForecast.Range("A2").Value = 45592.99999999999
$poForecast.Range("B2").Value = 20
$poForecast.Range("C2").Value = 19.99974458274913
$poForecast.Range("D2").Value = 19.99974459296693
$poForecast.Range("A3").Value = 45599.99999999999
$poForecast.Range("B3").Value = 220
$poForecast.Range("C3").Value = 219.999764614902
$poForecast.Range("D3").Value = 219.9997646254805
$poForecast.Range("A4").Value = 45606.99999999999
$poForecast.Range("B4").Value = 420
$poForecast.Range("C4").Value = 419.9997828128148
$poForecast.Range("D4").Value = 419.9997862926323
$poForecast.Range("A5").Value = 45613.99999999999
$poForecast.Range("B5").Value = 620
$poForecast.Range("C5").Value = 619.9997987350289
$poForecast.Range("D5").Value = 619.999810042814
$poForecast.Range("A6").Value = 45620.99999999999
$poForecast.Range("B6").Value = 820
$poForecast.Range("C6").Value = 819.9998136160418
$poForecast.Range("D6").Value = 819.9998356421606
$poForecast.Range("A7").Value = 45627.99999999999
$poForecast.Range("B7").Value = 1020
$poForecast.Range("C7").Value = 1019.999827094313
$poForecast.Range("D7").Value = 1019.999862175431
$poForecast.Range("A8").Value = 45634.99999999999
$poForecast.Range("B8").Value = 1220
$poForecast.Range("C8").Value = 1219.999839867692
$poForecast.Range("D8").Value = 1219.99989007196
$poForecast.Range("A9").Value = 45641.99999999999
$poForecast.Range("B9").Value = 1420
$poForecast.Range("C9").Value = 1419.999851602227
$poForecast.Range("D9").Value = 1419.999918338013
$poForecast.Range("A10").Value = 45648.99999999999
$poForecast.Range("B10").Value = 1620
$poForecast.Range("C10").Value = 1619.99986350101
$poForecast.Range("D10").Value = 1619.999946685036
$poForecast.Range("A11").Value = 45655.99999999999
$poForecast.Range("B11").Value = 1820
$poForecast.Range("C11").Value = 1819.999876703831
$poForecast.Range("D11").Value = 1819.999974871369

# Keep the original sheet active, like before the edit.
$weekly.Activate()
